$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 75, pushing existing rows 75-85 down to 77-87.
$ws.Rows.Item(75).Insert()
$ws.Rows.Item(75).Insert()

# New row 75: Santina / Especial, dated 44524
$ws.Range("A75").Value = 2
$ws.Range("B75").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C75").Value = "Coquimbo"
$ws.Range("D75").Value = 44524
$ws.Range("E75").Value = 4
$ws.Range("F75").Value = "Fruta"
$ws.Range("G75").Value = 100103
$ws.Range("H75").Value = "Frutos de hueso (carozo)"
$ws.Range("I75").Value = 100103001
$ws.Range("J75").Value = "Cereza"
$ws.Range("K75").Value = "Santina"
$ws.Range("L75").Value = "Especial"
$ws.Range("M75").Value = 400
$ws.Range("N75").Value = 23000
$ws.Range("O75").Value = 24000
$ws.Range("P75").Value = 23500
$ws.Range("Q75").Value = "$/bandeja 10 kilos"
$ws.Range("R75").Value = "Provincia de Curicó"
$ws.Range("S75").Value = 2350
$ws.Range("T75").Value = 10

# New row 76: Santina / Primera, dated 44524
$ws.Range("A76").Value = 2
$ws.Range("B76").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C76").Value = "Coquimbo"
$ws.Range("D76").Value = 44524
$ws.Range("E76").Value = 4
$ws.Range("F76").Value = "Fruta"
$ws.Range("G76").Value = 100103
$ws.Range("H76").Value = "Frutos de hueso (carozo)"
$ws.Range("I76").Value = 100103001
$ws.Range("J76").Value = "Cereza"
$ws.Range("K76").Value = "Santina"
$ws.Range("L76").Value = "Primera"
$ws.Range("M76").Value = 240
$ws.Range("N76").Value = 21000
$ws.Range("O76").Value = 22000
$ws.Range("P76").Value = 21500
$ws.Range("Q76").Value = "$/bandeja 10 kilos"
$ws.Range("R76").Value = "Provincia de Curicó"
$ws.Range("S76").Value = 2150
$ws.Range("T76").Value = 10
